$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Rows 12-15: four new "Exception" code entries (PLX0011 .. PLX0014)
# ----------------------------------------------------------------------

$ws.Range("B12").Value = "Exception"
$ws.Range("C12").Formula = "=_xlfn.SWITCH(B12,""Trace"",""T"",""Debug"",""D"",""Information"",""I"",""Warning"",""W"",""Error"",""E"",""Exception"",""X"")"
$ws.Range("D12").Formula = "=""PL""&C12&RIGHT(""0000""&A12,4)"
$ws.Range("E12").Value = "参数 “{0}”超出了指定的范围。"
$ws.Range("F12").Formula = "=D12&"": ""&E12"

$ws.Range("B13").Value = "Exception"
$ws.Range("C13").Formula = "=_xlfn.SWITCH(B13,""Trace"",""T"",""Debug"",""D"",""Information"",""I"",""Warning"",""W"",""Error"",""E"",""Exception"",""X"")"
$ws.Range("D13").Formula = "=""PL""&C13&RIGHT(""0000""&A13,4)"
$ws.Range("E13").Value = "参数 “{0}”超出了指定的范围 “{1}”。"
$ws.Range("F13").Formula = "=D13&"": ""&E13"

$ws.Range("B14").Value = "Exception"
$ws.Range("C14").Formula = "=_xlfn.SWITCH(B14,""Trace"",""T"",""Debug"",""D"",""Information"",""I"",""Warning"",""W"",""Error"",""E"",""Exception"",""X"")"
$ws.Range("D14").Formula = "=""PL""&C14&RIGHT(""0000""&A14,4)"
$ws.Range("E14").Value = "参数超出了指定的范围。"
$ws.Range("F14").Formula = "=D14&"": ""&E14"

$ws.Range("B15").Value = "Exception"
$ws.Range("C15").Formula = "=_xlfn.SWITCH(B15,""Trace"",""T"",""Debug"",""D"",""Information"",""I"",""Warning"",""W"",""Error"",""E"",""Exception"",""X"")"
$ws.Range("D15").Formula = "=""PL""&C15&RIGHT(""0000""&A15,4)"
$ws.Range("E15").Value = "参数超出了指定的范围 “{0}”。"
$ws.Range("F15").Formula = "=D15&"": ""&E15"

# ----------------------------------------------------------------------
# Data validation: extend the "Trace,Debug,..." drop-down list so it
# also covers the newly populated B12:B15 cells (in addition to the
# pre-existing B2:B11).
# ----------------------------------------------------------------------

$full = $ws.Range("B2:B15")
$full.Validation.Delete()
$full.Validation.Add(3, 1, 1, """Trace,Debug,Information,Warning,Error,Exception""")

Write-Host "done"
